# Updates the Neo4j/Cypher query text stored in column B (per-tab query) and
# column C (shared StatQuery) of the "startup" sheet, plus the row heights and
# active selection that change as a side effect of the longer query text.
#
# The query strings contain backticks, quotes, dollar signs and newlines, so
# to avoid any PowerShell quoting/escaping pitfalls each one is stored here as
# base64 and decoded at runtime.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Decode-B64([string]$s) {
    return [System.Text.Encoding]::UTF8.GetString([System.Convert]::FromBase64String($s))
}

# $q11 = new shared "StatQuery" (column C on every row)
$q11 = Decode-B64 'TUFUQ0ggKHA6cHJvZ3JhbSk8LS0oczpzdHVkeSk8LS0oYykKTUFUQ0ggKGNmKS0tPihzYW1wOnNhbXBsZSkKV0hFUkUgc2FtcC5zYW1wbGVfc2l0ZSBJTiBbJ0JsYWRkZXIgVHJpZ29uZSddCk1BVENIIChjZjpmaWxlKS1bKl0tPihjOmNhc2UpCk9QVElPTkFMIE1BVENIIChzZjpmaWxlKS0tPihzKQpSRVRVUk4KICAgIGNvdW50KGRpc3RpbmN0IHApIEFTIFByb2dyYW1zLAogICAgY291bnQoZGlzdGluY3QgcykgQVMgU3R1ZGllcywKICAgIGNvdW50KGRpc3RpbmN0IGMpIEFTIENhc2VzLAogICAgY291bnQoZGlzdGluY3Qgc2FtcCkgQVMgU2FtcGxlcywKICAgIGNvdW50KGRpc3RpbmN0IGNmKSBBUyBgQ2FzZSBGaWxlc2AsCiAgICBjb3VudChkaXN0aW5jdCBzZikgQVMgYFN0dWR5IEZpbGVzYA=='
# $q12 = new FilesTab query (row 4, column B)
$q12 = Decode-B64 'TUFUQ0ggKGY6ZmlsZSktLT4ocGFyZW50KQpXSVRIIERJU1RJTkNUIGYsIHBhcmVudApNQVRDSCAoZiktWypdLT4oYzpjYXNlKTwtLShkZW1vOmRlbW9ncmFwaGljKQpNQVRDSCAoczpzdHVkeSk8LVsqXS0oYyk8LS0oZGlhZzpkaWFnbm9zaXMpCk1BVENIIChmKS0tPihzYW1wOnNhbXBsZSkKV0hFUkUgc2FtcC5zYW1wbGVfc2l0ZSBJTiBbJ0JsYWRkZXIgVHJpZ29uZSddCiBNQVRDSCAoZiktWypdLT4oc2FtcDpzYW1wbGUpCldJVEgKICAgICAgICBESVNUSU5DVCBmLCBwYXJlbnQsIGMsIGRlbW8sIGRpYWcsIHMsIHNhbXAsCiAgICAgICAgWydCeXRlcycsICdLQicsICdNQicsICdHQicsICdUQiddIEFTIHVuaXRzLAogICAgICAgIHRvSW50ZWdlcihmbG9vcihsb2coZi5maWxlX3NpemUpL2xvZygxMDI0KSkpIGFzIGksCiAgICAgICAgMiBhcyBwcmVjaXNpb24KV0lUSAogICAgICAgIGYsIHBhcmVudCwgYywgZGVtbywgZGlhZywgcywgc2FtcCwKICAgICAgICBmLmZpbGVfc2l6ZSAvKDEwMjReaSkgQVMgdmFsdWUsCiAgICAgICAgMTBecHJlY2lzaW9uIEFTIGZhY3RvciwKICAgICAgICB1bml0c1tpXSBhcyB1bml0CldJVEgKICAgICAgICBmLCBwYXJlbnQsIGMsIGRlbW8sIGRpYWcsIHMsIHNhbXAsIHVuaXQsCiAgICAgICAgcm91bmQoZmFjdG9yICogdmFsdWUpL2ZhY3RvciBBUyBzaXplClJFVFVSTgogICAgICAgIGNvYWxlc2NlKGYuZmlsZV9uYW1lLCAnJykgQVMgYEZpbGUgTmFtZWAsCiAgICAgICAgY29hbGVzY2UoZi5maWxlX2Zvcm1hdCwgJycpIEFTIGBGb3JtYXRgLAogICAgICAgIGNvYWxlc2NlKGYuZmlsZV90eXBlLCAnJykgQVMgYEZpbGUgVHlwZWAsCiAgICAgICAgQ0FTRSBzaXplICUgMSBXSEVOIDAgVEhFTiBhcG9jLmNvbnZlcnQudG9JbnRlZ2VyKHNpemUpKycgJyArdW5pdCBFTFNFIHNpemUrJyAnICt1bml0IEVORCBBUyBTaXplLAogICAgICAgIGNvYWxlc2NlKGxhYmVscyhwYXJlbnQpWzBdLCAnJykgQVMgYEFzc29jaWF0aW9uYCwKICAgICAgICBjb2FsZXNjZShmLmZpbGVfZGVzY3JpcHRpb24sICcnKSBBUyBgRGVzY3JpcHRpb25gLAogICAgICAgIGNvYWxlc2NlKHNhbXAuc2FtcGxlX2lkLCAnJykgQVMgYFNhbXBsZSBJRGAsCiAgICAgICAgY29hbGVzY2UoYy5jYXNlX2lkLCAnJykgQVMgYENhc2UgSURgLAogICAgICAgIGNvYWxlc2NlKGRlbW8uYnJlZWQsJycpIEFTIEJyZWVkICwKICAgICAgICBjb2FsZXNjZShkaWFnLmRpc2Vhc2VfdGVybSwnJykgQVMgRGlhZ25vc2lzCiAgICAgICAgb3JkZXIgYnkgZi5maWxlX25hbWUgYXNjCiAgICAgICAgbGltaXQgMTAw'
# $q13 = updated SamplesTab query (row 3, column B)
$q13 = Decode-B64 'TUFUQ0ggKHM6c3R1ZHkpPC1bKl0tKGM6Y2FzZSk8LS0oZGVtbzpkZW1vZ3JhcGhpYyksIChzYW1wOnNhbXBsZSktLT4oYyk8LS0oZGlhZzpkaWFnbm9zaXMpIApXSEVSRSAgc2FtcC5zYW1wbGVfc2l0ZSBJTiBbJ0JsYWRkZXIgVHJpZ29uZSddCldJVEggRElTVElOQ1Qgc2FtcCBBUyBzYW1wLCBjLCBkZW1vLCBkaWFnClJFVFVSTiAgY29hbGVzY2Uoc2FtcC5zYW1wbGVfaWQsICcnKSBBUyBgU2FtcGxlIElEYCwgCiAgICAgICAgY29hbGVzY2UoYy5jYXNlX2lkLCAnJykgQVMgYENhc2UgSURgLCAKICAgICAgICBjb2FsZXNjZShkZW1vLmJyZWVkLCcnKSBBUyBCcmVlZCwKICAgICAgICBjb2FsZXNjZShkaWFnLmRpc2Vhc2VfdGVybSwnJykgQVMgRGlhZ25vc2lzLCAKICAgICAgICBjb2FsZXNjZShzYW1wLnNhbXBsZV9zaXRlLCAnJykgQVMgYFNhbXBsZSBTaXRlYCwKICAgICAgICBjb2FsZXNjZShzYW1wLnN1bW1hcml6ZWRfc2FtcGxlX3R5cGUsICcnKSBBUyBgU2FtcGxlIFR5cGVgLAogICAgICAgIGNvYWxlc2NlKHNhbXAuc3BlY2lmaWNfc2FtcGxlX3BhdGhvbG9neSwgJycpIEFTIGBQYXRob2xvZ3kvTW9ycGhvbG9neWAsCiAgICAgICAgY29hbGVzY2Uoc2FtcC50dW1vcl9ncmFkZSwgJycpIEFTIGBUdW1vciBHcmFkZWAsCiAgICAgICAgY29hbGVzY2Uoc2FtcC5zYW1wbGVfY2hyb25vbG9neSwgJycpIEFTIGBTYW1wbGUgQ2hyb25vbG9neWAsCiAgICAgICAgY29hbGVzY2Uoc2FtcC5wZXJjZW50YWdlX3R1bW9yLCAnJykgQVMgYFBlcmNlbnRhZ2UgVHVtb3JgLAogICAgICAgIGNvYWxlc2NlKHNhbXAubmVjcm9wc3lfc2FtcGxlLCAnJykgQVMgYE5lY3JvcHN5IFNhbXBsZWAsCiAgICAgICAgY29hbGVzY2Uoc2FtcC5zYW1wbGVfcHJlc2VydmF0aW9uLCAnJykgQVMgYFNhbXBsZSBQcmVzZXJ2YXRpb25gCm9yZGVyIGJ5IHNhbXAuc2FtcGxlX2lkIGFzYwpsaW1pdCAxMDA='
# $q14 = updated StudyFilesTab query (row 5, column B)
$q14 = Decode-B64 'TUFUQ0ggKGY6ZmlsZSktLT4oczpzdHVkeSkKTUFUQ0ggKHMpPC0tKGM6Y2FzZSk8LS0oZGlhZzpkaWFnbm9zaXMpCk1BVENIIChjKTwtLShkZW1vOmRlbW9ncmFwaGljKQpNQVRDSCAoc2FtcDpzYW1wbGUpLS0+KGMpCldIRVJFIHNhbXAuc2FtcGxlX3NpdGUgSU4gWydCbGFkZGVyIFRyaWdvbmUnXQpXSVRICiAgICAgICAgRElTVElOQ1QgZiwgYywgZGVtbywgZGlhZywgcywKICAgICAgICBbJ0J5dGVzJywgJ0tCJywgJ01CJywgJ0dCJywgJ1RCJ10gQVMgdW5pdHMsCiAgICAgICAgdG9JbnRlZ2VyKGZsb29yKGxvZyhmLmZpbGVfc2l6ZSkvbG9nKDEwMjQpKSkgYXMgaSwKICAgICAgICAyIGFzIHByZWNpc2lvbgpXSVRICiAgICAgICAgZiwgYywgZGVtbywgZGlhZywgcywKICAgICAgICBmLmZpbGVfc2l6ZSAvKDEwMjReaSkgQVMgdmFsdWUsIDEwXnByZWNpc2lvbiBBUyBmYWN0b3IsCiAgICAgICAgdW5pdHNbaV0gYXMgdW5pdAogICAgICAgIFdJVEgKICAgICAgICBmLCAgYywgZGVtbywgZGlhZywgcywgdW5pdCwKICAgICAgICByb3VuZChmYWN0b3IgKiB2YWx1ZSkvZmFjdG9yIEFTIHNpemUKUkVUVVJOIERJU1RJTkNUCiAgY29hbGVzY2UoZi5maWxlX25hbWUsICcnKSBBUyBgRmlsZSBOYW1lYCwKICBjb2FsZXNjZShmLmZpbGVfdHlwZSwgJycpIEFTIGBGaWxlIFR5cGVgLAogIGNvYWxlc2NlKCJzdHVkeSIsICcnKSBBUyBgQXNzb2NpYXRpb25gLAogIGNvYWxlc2NlKGYuZmlsZV9kZXNjcmlwdGlvbiwgJycpIEFTIGBEZXNjcmlwdGlvbmAsCiAgY29hbGVzY2UoZi5maWxlX2Zvcm1hdCwgJycpIEFTICBGb3JtYXQsCiAgQ0FTRSBzaXplICUgMSBXSEVOIDAgVEhFTiBhcG9jLmNvbnZlcnQudG9JbnRlZ2VyKHNpemUpKycgJyArdW5pdCBFTFNFIHNpemUrJyAnICt1bml0IEVORCBBUyBTaXplLAogIGNvYWxlc2NlKHMuY2xpbmljYWxfc3R1ZHlfZGVzaWduYXRpb24sJycpIEFTIGBTdHVkeSBDb2RlYAogIG9yZGVyIGJ5ICdGaWxlIE5hbWUnIGFzYwogIGxpbWl0IDEwMA=='
# $q15 = new CasesTab query (row 2, column B)
$q15 = Decode-B64 'TUFUQ0ggKHM6c3R1ZHkpPC1bKl0tKGM6Y2FzZSk8LS0oZGVtbzpkZW1vZ3JhcGhpYykKTUFUQ0ggKGMpPC0tKGRpYWc6ZGlhZ25vc2lzKQpPUFRJT05BTCBNQVRDSCAoY286Y29ob3J0KTwtWypdLShjKQpNQVRDSCAoc2FtcDpzYW1wbGUpLS0+KGMpCldIRVJFICBzYW1wLnNhbXBsZV9zaXRlIGluIFsnQmxhZGRlciBUcmlnb25lJ10KV0lUSCBESVNUSU5DVCBjLCBzLCBkZW1vLCBkaWFnLCBjbywgZGVtby5wYXRpZW50X2FnZV9hdF9lbnJvbGxtZW50IEFTIGFnZSwgZGVtby53ZWlnaHQgYXMgd2VpZ2h0ClJFVFVSTiAgY29hbGVzY2UoYy5jYXNlX2lkLCAnJykgQVMgYENhc2UgSURgICwKICAgICAgICBjb2FsZXNjZShzLmNsaW5pY2FsX3N0dWR5X2Rlc2lnbmF0aW9uLCAnJykgQVMgYFN0dWR5IENvZGVgICwKICAgICAgICBjb2FsZXNjZShzLmNsaW5pY2FsX3N0dWR5X3R5cGUsICcnKSBBUyAgYFN0dWR5IFR5cGVgLAogICAgICAgIGNvYWxlc2NlKGRlbW8uYnJlZWQsICcnKSBBUyBCcmVlZCAsCiAgICAgICAgY29hbGVzY2UoZGlhZy5kaXNlYXNlX3Rlcm0sICcnKSBBUyBEaWFnbm9zaXMgLAogICAgICAgIGNvYWxlc2NlKGRpYWcuc3RhZ2Vfb2ZfZGlzZWFzZSwgJycpIEFTIGBTdGFnZSBvZiBEaXNlYXNlYCAsCiAgY29hbGVzY2UoQ0FTRSBhZ2UgJSAxIFdIRU4gMCBUSEVOIGFwb2MuY29udmVydC50b0ludGVnZXIoYWdlKSBFTFNFIGFnZSBFTkQsICcnKSBBUyBBZ2UsCiAgICAgICBjb2FsZXNjZShkZW1vLnNleCwgJycpIEFTIFNleCwKICAgICAgIGNvYWxlc2NlKGRlbW8ubmV1dGVyZWRfaW5kaWNhdG9yLCAnJykgQVMgYE5ldXRlcmVkIFN0YXR1c2AsCmNvYWxlc2NlKENBU0Ugd2VpZ2h0ICUgMSBXSEVOIDAgVEhFTiBhcG9jLmNvbnZlcnQudG9JbnRlZ2VyKHdlaWdodCkgRUxTRSB3ZWlnaHQgRU5ELCAnJykgQVMgYFdlaWdodCAoa2cpYCwKICAgICAgIGNvYWxlc2NlKGRpYWcuYmVzdF9yZXNwb25zZSwgJycpIEFTIGBSZXNwb25zZSB0byBUcmVhdG1lbnRgLAogICAgICAgY29hbGVzY2UoY28uY29ob3J0X2Rlc2NyaXB0aW9uLCAnJykgQVMgYENvaG9ydGAKb3JkZXIgYnkgYy5jYXNlX2lkIGFzYwpsaW1pdCAxMDA='

# Row 2 (CasesTab): B2 -> new Cases query (was si15), C2 -> new StatQuery (was si11)
$ws.Range("B2").Value = $q15
$ws.Range("C2").Value = $q11

# Row 3 (SamplesTab): B3 -> updated Samples query, C3 -> new StatQuery
$ws.Range("B3").Value = $q13
$ws.Range("C3").Value = $q11

# Row 4 (FilesTab): B4 -> new Files query (was si12 StatQuery slot), C4 -> new StatQuery
$ws.Range("B4").Value = $q12
$ws.Range("C4").Value = $q11

# Row 5 (StudyFilesTab): B5 -> updated StudyFiles query, C5 -> new StatQuery
$ws.Range("B5").Value = $q14
$ws.Range("C5").Value = $q11

# Row heights
$ws.Rows.Item(2).RowHeight = 333.5
$ws.Rows.Item(3).RowHeight = 275.5
$ws.Rows.Item(4).RowHeight = 409.5
$ws.Rows.Item(5).RowHeight = 409.5

# Update selection / view (scrolls back up, selects C2)
[void]$ws.Range("C2").Select()

Write-Host "edit applied"
